$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.405.05"
$ws.Range("E2").Value = "  +0.48%  "

# Row 3
$ws.Range("D3").Value = "1.869.44"
$ws.Range("E3").Value = "  +0.20%  "

# Row 4
$ws.Range("D4").Value = "'1.001"

# Row 5
$ws.Range("D5").Value = "'246.99"
$ws.Range("E5").Value = "  +1.86%  "

# Row 6
$ws.Range("D6").Value = "'1.000"

# Row 7
$ws.Range("D7").Value = "'0.4725"
$ws.Range("E7").Value = "  +0.17%  "

# Row 8
$ws.Range("D8").Value = "'0.2913"
$ws.Range("E8").Value = "  +1.64%  "

# Row 9
$ws.Range("D9").Value = "'0.06480"
$ws.Range("E9").Value = "  +0.21%  "

# Row 10
$ws.Range("D10").Value = "'22.05"
$ws.Range("E10").Value = "  +6.03%  "

# Row 11
$ws.Range("D11").Value = "'0.07722"
$ws.Range("E11").Value = "  +0.05%  "

# Row 12
$ws.Range("D12").Value = "'97.52"
$ws.Range("E12").Value = "  +2.73%  "

# Row 13
$ws.Range("D13").Value = "'0.7405"
$ws.Range("E13").Value = "  +4.88%  "

# Row 14
$ws.Range("D14").Value = "1.870.09"
$ws.Range("E14").Value = "  +0.18%  "

# Row 15
$ws.Range("D15").Value = "'5.141"
$ws.Range("E15").Value = "  +1.16%  "

# Row 16
$ws.Range("D16").Value = "'272.86"
$ws.Range("E16").Value = "  +1.29%  "

# Row 17
$ws.Range("D17").Value = "30.394.96"
$ws.Range("E17").Value = "  +0.46%  "

# Row 18
$ws.Range("D18").Value = "'13.38"
$ws.Range("E18").Value = "  +0.50%  "

# Row 19
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.07%  "

# Row 20
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.000007501"
$ws.Range("E20").Value = "  -0.46%  "

# Row 21
$ws.Range("D21").Value = "2.118.92"
$ws.Range("E21").Value = "  +0.11%  "

# Row 22
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.14%  "

# Row 23
$ws.Range("D23").Value = "'5.236"
$ws.Range("E23").Value = "  +0.90%  "

# Row 24
$ws.Range("D24").Value = "'6.169"
$ws.Range("E24").Value = "  +0.83%  "

# Row 25
$ws.Range("D25").Value = "'9.263"
$ws.Range("E25").Value = "  -0.66%  "

# Row 26
$ws.Range("D26").Value = "'163.36"
$ws.Range("E26").Value = "  -1.14%  "

# Row 27
$ws.Range("D27").Value = "'18.74"
$ws.Range("E27").Value = "  -0.46%  "

# Row 28
$ws.Range("D28").Value = "'1.919"
$ws.Range("E28").Value = "  +0.46%  "

# Row 29
$ws.Range("D29").Value = "'0.1001"
$ws.Range("E29").Value = "  +1.67%  "

# Row 30
$ws.Range("D30").Value = "'1.371"
$ws.Range("E30").Value = "  -0.47%  "

# Row 31
$ws.Range("D31").Value = "'1.502"
$ws.Range("E31").Value = "  -0.22%  "

# Row 32
$ws.Range("D32").Value = "'4.269"
$ws.Range("E32").Value = "  +0.62%  "

# Row 33
$ws.Range("D33").Value = "'4.116"
$ws.Range("E33").Value = "  +2.56%  "

# Row 34
$ws.Range("D34").Value = "'0.04813"
$ws.Range("E34").Value = "  +1.54%  "

# Row 35
$ws.Range("D35").Value = "'1.116"
$ws.Range("E35").Value = "  -0.15%  "

# Row 36
$ws.Range("D36").Value = "'0.6939"
$ws.Range("E36").Value = "  +0.52%  "

# Row 37
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.714"
$ws.Range("E37").Value = "  +0.35%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01846"
$ws.Range("E38").Value = "  +0.18%  "

# Row 39
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.742"
$ws.Range("E39").Value = "  +0.36%  "

# Row 40
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'6.275"
$ws.Range("E40").Value = "  -0.83%  "

# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'72.80"
$ws.Range("E41").Value = "  +3.56%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'1.966"
$ws.Range("E42").Value = "  +3.92%  "

# Row 43
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4177"
$ws.Range("E43").Value = "  +2.75%  "

# Row 44
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  +0.06%  "

# Row 45
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'0.8364"
$ws.Range("E45").Value = "  -0.38%  "

# Row 46
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'101.87"
$ws.Range("E46").Value = "  -0.37%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'9.346"
$ws.Range("E47").Value = "  +1.42%  "

# Row 48
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "'35.41"
$ws.Range("E48").Value = "  +2.05%  "

# Row 49
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'6.964"
$ws.Range("E49").Value = "  -1.37%  "

# Row 50
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "'921.11"
$ws.Range("E50").Value = "  -0.89%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05639"
$ws.Range("E51").Value = "  +1.36%  "
